$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    -0.26275633660550568,
    -0.17764675689583775,
    -0.072934656668248721,
    -0.13452930694590037,
    -0.12783831554415315,
    -0.0018633111805912428,
    0.018141837413656603,
    0.038200044582343828,
    0.044409664309439378,
    0.050667317382405486,
    0.055222429821775165,
    -0.010022512059513655,
    -0.0040127816475434841,
    0.0079884476250455805,
    -0.021050216573242686,
    -0.015026194408048976,
    -0.009004304729754331,
    -0.074443307220576571,
    -0.0653602353520637,
    -0.018014158077413356,
    -0.0090043660781677914,
    -0.093946223948108809,
    -0.084633472576572721,
    -0.042126678116074068,
    -0.11860210209372468,
    -0.11223815074758647,
    -0.082370779569619135,
    -0.075291479150605412,
    -0.06254356497444391,
    -0.042172456125110003,
    -0.02702041050634385,
    -0.0060007833746036354
)

$colB = @(
    0.26209006382262601,
    0.17588088032066729,
    0.072532645400791296,
    0.13383831525510814,
    0.12646541297734792,
    0.0018581622343569748,
    -0.018200044933694315,
    -0.038409664603356042,
    -0.044667317673400042,
    -0.050722430106496574,
    -0.055349553414210106,
    0.010012781357204403,
    0.0040115520609038668,
    -0.0079970203131933815,
    0.02102619411688611,
    0.015004304437423954,
    0.0089999996956278139,
    0.074360235061568858,
    0.064735483373748792,
    0.018004365780686982,
    0.0089999997022234268,
    0.093633472278217056,
    0.084126677683371121,
    0.041999999564986013,
    0.11823815044909836,
    0.11176337209074916,
    0.081291478843177334,
    0.074543564638664961,
    0.06217245575560959,
    0.04202041015336988,
    0.027000782997557593,
    0.0059999996817685641
)

for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}
